$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new data rows right before the current row 56. This pushes the
# existing rows 56..68 down to 58..70 (carrying their values/formatting with
# them), and leaves two blank rows at 56 and 57 for the new data below.
$ws.Rows("56:57").Insert()

# --- Row 56 (brand new record) ---
$ws.Range("A56").Value = 8
$ws.Range("B56").Value = "Terminal La Palmera de La Serena"
$ws.Range("C56").Value = "Coquimbo"
$ws.Range("D56").Value = 44900
$ws.Range("E56").Value = 4
$ws.Range("F56").Value = "Fruta"
$ws.Range("G56").Value = 100103
$ws.Range("H56").Value = "Frutos de hueso (carozo)"
$ws.Range("I56").Value = 100103003
$ws.Range("J56").Value = "Damasco"
$ws.Range("K56").Value = "Castle Brite"
$ws.Range("L56").Value = "Primera"
$ws.Range("M56").Value = 300
$ws.Range("N56").Value = 24000
$ws.Range("O56").Value = 25000
$ws.Range("P56").Value = 24500
$ws.Range("Q56").Value = "$/caja 18 kilos"
$ws.Range("R56").Value = "Provincia de San Felipe de Aconcagua"
$ws.Range("S56").Value = 1361
$ws.Range("T56").Value = 18

# --- Row 57 (brand new record) ---
$ws.Range("A57").Value = 8
$ws.Range("B57").Value = "Terminal La Palmera de La Serena"
$ws.Range("C57").Value = "Coquimbo"
$ws.Range("D57").Value = 44900
$ws.Range("E57").Value = 4
$ws.Range("F57").Value = "Fruta"
$ws.Range("G57").Value = 100103
$ws.Range("H57").Value = "Frutos de hueso (carozo)"
$ws.Range("I57").Value = 100103003
$ws.Range("J57").Value = "Damasco"
$ws.Range("K57").Value = "Castle Brite"
$ws.Range("L57").Value = "Segunda"
$ws.Range("M57").Value = 280
$ws.Range("N57").Value = 19000
$ws.Range("O57").Value = 20000
$ws.Range("P57").Value = 19500
$ws.Range("Q57").Value = "$/caja 18 kilos"
$ws.Range("R57").Value = "Provincia de San Felipe de Aconcagua"
$ws.Range("S57").Value = 1083
$ws.Range("T57").Value = 18
